# Applies the odds/score updates from the 2025-04-16 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8
$ws.Range("J8").Value = 1.04
$ws.Range("K8").Value = 13
$ws.Range("L8").Value = 1.22
$ws.Range("M8").Value = 4.33
$ws.Range("N8").Value = 1.67
$ws.Range("O8").Value = 2.15
$ws.Range("P8").Value = 1.33
$ws.Range("Q8").Value = 3.25
$ws.Range("R8").Value = 2.5
$ws.Range("S8").Value = 1.5
$ws.Range("T8").Value = 6.5
$ws.Range("Z8").Value = 12
$ws.Range("AE8").Value = 23

# Row 9
$ws.Range("G9").Value = 2.5
$ws.Range("I9").Value = 3.1
$ws.Range("L9").Value = 1.5
$ws.Range("M9").Value = 2.63
$ws.Range("N9").Value = 2.4
$ws.Range("O9").Value = 1.53
$ws.Range("R9").Value = 2.05
$ws.Range("S9").Value = 1.7
$ws.Range("Y9").Value = 41
$ws.Range("Z9").Value = 7
$ws.Range("AD9").Value = 501
$ws.Range("AG9").Value = 12

# Row 10
$ws.Range("G10").Value = 2.3
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 3.4
$ws.Range("L10").Value = 1.53
$ws.Range("M10").Value = 2.5
$ws.Range("U10").Value = 9.5
$ws.Range("W10").Value = 21
$ws.Range("AE10").Value = 8
$ws.Range("AG10").Value = 13
$ws.Range("AH10").Value = 41
$ws.Range("AI10").Value = 34

# Row 13
$ws.Range("I13").Value = 3.1

# Row 22
$ws.Range("G22").Value = 1.75
$ws.Range("H22").Value = 3.4
$ws.Range("I22").Value = 4.75
$ws.Range("K22").Value = 7.5
$ws.Range("L22").Value = 1.4
$ws.Range("M22").Value = 2.75
$ws.Range("P22").Value = 1.5
$ws.Range("Q22").Value = 2.5
$ws.Range("R22").Value = 2.2
$ws.Range("S22").Value = 1.62
$ws.Range("T22").Value = 5.5
$ws.Range("U22").Value = 7.5
$ws.Range("V22").Value = 9
$ws.Range("W22").Value = 13
$ws.Range("X22").Value = 17
$ws.Range("AA22").Value = 7
$ws.Range("AB22").Value = 21
$ws.Range("AC22").Value = 81
$ws.Range("AE22").Value = 10
$ws.Range("AF22").Value = 23
$ws.Range("AG22").Value = 17
$ws.Range("AH22").Value = 51
$ws.Range("AJ22").Value = 51

# Row 25
$ws.Range("G25").Value = 3.25
$ws.Range("H25").Value = 3.3
$ws.Range("I25").Value = 2.1
$ws.Range("K25").Value = 9.5
$ws.Range("L25").Value = 1.3
$ws.Range("M25").Value = 3.4
$ws.Range("P25").Value = 1.4
$ws.Range("Q25").Value = 2.75
$ws.Range("R25").Value = 1.8
$ws.Range("S25").Value = 1.91
$ws.Range("X25").Value = 26
$ws.Range("Z25").Value = 9.5
$ws.Range("AD25").Value = 251
$ws.Range("AI25").Value = 19

# Row 26
$ws.Range("G26").Value = 1.38
$ws.Range("H26").Value = 4.75
$ws.Range("I26").Value = 6.5
$ws.Range("J26").Value = 1.03
$ws.Range("L26").Value = 1.22
$ws.Range("R26").Value = 2.1
$ws.Range("S26").Value = 1.67
$ws.Range("U26").Value = 6.5
$ws.Range("W26").Value = 9
$ws.Range("AA26").Value = 9.5
$ws.Range("AC26").Value = 67
$ws.Range("AD26").Value = 1000
$ws.Range("AE26").Value = 17
$ws.Range("AF26").Value = 34

# Row 41
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 3.25
$ws.Range("I41").Value = 2.45
$ws.Range("N41").Value = 1.85
$ws.Range("O41").Value = 2
$ws.Range("P41").Value = 1.36
$ws.Range("Q41").Value = 3
$ws.Range("R41").Value = 1.62
$ws.Range("S41").Value = 2.2
$ws.Range("T41").Value = 11
$ws.Range("X41").Value = 23
$ws.Range("AG41").Value = 9.5
$ws.Range("AH41").Value = 23

# Row 42
$ws.Range("T42").Value = 9.5

# Row 43
$ws.Range("I43").Value = 1.67
$ws.Range("J43").Value = 1.01
$ws.Range("L43").Value = 1.08
$ws.Range("M43").Value = 6.5
$ws.Range("N43").Value = 1.4
$ws.Range("O43").Value = 2.88
$ws.Range("R43").Value = 1.44
$ws.Range("S43").Value = 2.63
$ws.Range("AC43").Value = 29
$ws.Range("AD43").Value = 81
$ws.Range("AG43").Value = 9

Write-Host "Applied 113 cell updates"
